$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 54, shifting existing rows 54+ down to 56+.
$ws.Rows("54:55").Insert()

# Fixed columns shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"

# New row 54: Hass / Primera, Perú, bandeja 10 kilos
$ws.Cells.Item(54, 1).Value  = $mercadoId
$ws.Cells.Item(54, 2).Value  = $mercado
$ws.Cells.Item(54, 3).Value  = $region
$ws.Cells.Item(54, 4).Value  = 44771
$ws.Cells.Item(54, 5).Value  = $codreg
$ws.Cells.Item(54, 6).Value  = $tipo
$ws.Cells.Item(54, 7).Value  = $productoId
$ws.Cells.Item(54, 8).Value  = $producto
$ws.Cells.Item(54, 9).Value  = $categoriaId
$ws.Cells.Item(54, 10).Value = $categoria
$ws.Cells.Item(54, 11).Value = "Hass"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 200
$ws.Cells.Item(54, 14).Value = 27000
$ws.Cells.Item(54, 15).Value = 28000
$ws.Cells.Item(54, 16).Value = 27500
$ws.Cells.Item(54, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(54, 18).Value = "Perú"
$ws.Cells.Item(54, 19).Value = 2750
$ws.Cells.Item(54, 20).Value = 10

# New row 55: Hass / Segunda, Perú, bandeja 10 kilos
$ws.Cells.Item(55, 1).Value  = $mercadoId
$ws.Cells.Item(55, 2).Value  = $mercado
$ws.Cells.Item(55, 3).Value  = $region
$ws.Cells.Item(55, 4).Value  = 44771
$ws.Cells.Item(55, 5).Value  = $codreg
$ws.Cells.Item(55, 6).Value  = $tipo
$ws.Cells.Item(55, 7).Value  = $productoId
$ws.Cells.Item(55, 8).Value  = $producto
$ws.Cells.Item(55, 9).Value  = $categoriaId
$ws.Cells.Item(55, 10).Value = $categoria
$ws.Cells.Item(55, 11).Value = "Hass"
$ws.Cells.Item(55, 12).Value = "Segunda"
$ws.Cells.Item(55, 13).Value = 400
$ws.Cells.Item(55, 14).Value = 22000
$ws.Cells.Item(55, 15).Value = 23000
$ws.Cells.Item(55, 16).Value = 22500
$ws.Cells.Item(55, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(55, 18).Value = "Perú"
$ws.Cells.Item(55, 19).Value = 2250
$ws.Cells.Item(55, 20).Value = 10

# Match the date-format style used by the rest of column D.
$ws.Range("D54:D55").NumberFormat = $ws.Range("D56").NumberFormat
